# Populate ScenarioResults (sheet1) with scenario rows 2-16
$wb = $excel.ActiveWorkbook
$wsScenario = $wb.Worksheets.Item("ScenarioResults")

$scenarioData = @(
    @(1, 'https://dada.com.ua/', 'The operation wasn''t requested', 'The operation wasn''t requested', 'Id ISN''T on the page', 'Tag IS on the page', 'test'),
    @(2, 'https://elit-ampir.com.ua/ua/p1212509987-molding-home-decor.html', 'The operation wasn''t requested', 'The operation wasn''t requested', 'Id ISN''T on the page', 'Tag IS on the page', 'test'),
    @(3, 'https://www.youtube.com/watch?v=HfBJ0_1c4PU', 'The operation wasn''t requested', 'The operation wasn''t requested', 'The operation wasn''t requested', 'Tag IS on the page', 'test'),
    @(4, 'https://starkovden.github.io/Use-GitHub-Desktop.html', 'The operation wasn''t requested', 'The operation wasn''t requested', 'Id IS on the page', 'Tag IS on the page', 'test'),
    @(5, 'https://ru.wikipedia.org/wiki/The_Game_Awards', 'Uncorrectly requests', 'The operation wasn''t requested', 'The operation wasn''t requested', 'Tag IS on the page', 'test'),
    @(6, 'https://www.freecodecamp.org/ukrainian/news/yak-stvoryty-telehram-bota-za-dopomohoyu-python/', 'The operation wasn''t requested', 'The operation wasn''t requested', 'Id IS on the page', 'Tag IS on the page', 'test'),
    @(7, 'https://disted.edu.vn.ua/courses/learn/12688', 'The operation wasn''t requested', 'The operation wasn''t requested', 'The operation wasn''t requested', 'The operation wasn''t requested', 'test'),
    @(8, 'https://disted.edu.vn.ua/courses/learn/12688', 'The operation wasn''t requested', 'The operation wasn''t requested', 'Id IS on the page', 'Tag IS on the page', 'test'),
    @(9, 'https://habr.com/ru/articles/710350/', 'The operation wasn''t requested', 'The operation wasn''t requested', 'Id ISN''T on the page', 'Tag IS on the page', 'test'),
    @(10, 'https://github.com/Okarpets/NG_2023_Oleksandr_Karpets', 'The operation wasn''t requested', 'Uncorrectly requests', 'The operation wasn''t requested', 'Tag ISN''T on the page', 'test'),
    @(11, 'https://ru.stackoverflow.com/questions/1157100/post-%D0%B7%D0%B0%D0%BF%D1%80%D0%BE%D1%81-%D0%BA-api-youtube-%D0%BD%D0%B0-java', 'The operation wasn''t requested', 'The operation wasn''t requested', 'The operation wasn''t requested', 'Tag IS on the page', 'test'),
    @(12, 'https://ru.stackoverflow.com/questions/1159043/%d0%9f%d1%80%d0%be%d1%81%d1%82%d0%b5%d0%b9%d1%88%d0%b8%d0%b9-post-%d0%b7%d0%b0%d0%bf%d1%80%d0%be%d1%81-java-android?rq=1', 'The operation wasn''t requested', 'The operation wasn''t requested', 'Id IS on the page', 'The operation wasn''t requested', 'test'),
    @(13, 'https://www.rockstargames.com/ru/reddeadredemption2', 'The operation wasn''t requested', 'The operation wasn''t requested', 'Id ISN''T on the page', 'Tag ISN''T on the page', 'test'),
    @(14, 'https://ravesli.com/uroki-cpp/', 'The operation wasn''t requested', 'The operation wasn''t requested', 'Id IS on the page', 'Tag IS on the page', 'test'),
    @(15, 'https://www.microsoft.com/ru-ru/sql-server/sql-server-downloads', 'The operation wasn''t requested', 'The operation wasn''t requested', 'Id IS on the page', 'Tag IS on the page', 'test'),
)

for ($i = 0; $i -lt $scenarioData.Count; $i++) {
    $rowIndex = $i + 2
    $rowValues = $scenarioData[$i]
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $wsScenario.Cells.Item($rowIndex, $j + 1).Value = $rowValues[$j]
    }
}

# Populate CodePostResults (sheet "CodePostResults") with two 403 rows
$wsPost = $wb.Worksheets.Item("CodePostResults")
$wsPost.Cells.Item(2, 1).Value = 'https://github.com/Okarpets/NG_2023_Oleksandr_Karpets'
$wsPost.Cells.Item(2, 2).Value = 403
$wsPost.Cells.Item(3, 1).Value = 'https://github.com/Okarpets/NG_2023_Oleksandr_Karpets'
$wsPost.Cells.Item(3, 2).Value = 403
